# Add sector 0 = consumption rows to the raw sectors data sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows: final_sector(A)=0, NAICS(B)=blank, WIOD_sector(C)=..., name(D)=...
$newRows = @(
    @{ Row = 92; C = 37; D = "Final consumption expenditure by households" }
    @{ Row = 93; C = 38; D = "Final consumption expenditure by non-profit organisations serving households (NPISH)" }
    @{ Row = 94; C = 39; D = "Final consumption expenditure by government" }
    @{ Row = 95; C = 41; D = "Gross fixed capital formation" }
    @{ Row = 96; C = 42; D = "Changes in inventories and valuables" }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = 0
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}

# Match the scrolled/selected view state recorded in the saved file
$ws.Application.ActiveWindow.ScrollRow = 84
$ws.Range("B92").Select()
